$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")

# Row 31 (new fm26 test) needs to be populated first so that the original
# shared string (index 112, "Direct level from fm24 testing % TIV
# deductibles") is mutated in place into the "calcrule 6" variant and
# picked up by C31 (it is the sole referrer of that string before this
# edit). Doing C31 before C30 reproduces the shared-string layout in the
# target workbook (new string appended for C30's updated text instead).
$ws.Range("C31").Value = "Direct level from fm24 testing % TIV deductibles calcrule 6"
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = 6
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = "in progress"
$ws.Range("I31").Value = "in progress"

# Row 30 (fm25) description is now more specific about limits, which pushes
# its string into a newly appended shared-string slot.
$ws.Range("C30").Value = "Direct level from fm24 testing % TIV deductibles with limits calcrule 4"

# Selection moves to C32 like in the saved workbook.
$ws.Range("C32").Select
